$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix organization of nominal structures spend numbering in column E
$ws.Range("E38").Value = 9912
$ws.Range("E39").Value = 99121
$ws.Range("E40").Value = 991211
$ws.Range("E41").Value = 991212
$ws.Range("E42").Value = 99122
$ws.Range("E43").Value = 991221
$ws.Range("E44").Value = 991222
$ws.Range("E45").Value = 991223
$ws.Range("E46").Value = 991224
$ws.Range("E47").Value = 991225
$ws.Range("E48").Value = 9999
$ws.Range("E49").Value = 99991
$ws.Range("E50").Value = 99992

# Update the active selection to match the new cursor position after edits
$ws.Range("E48").Select()
